# Update Name of Algo
# Apply updated values to column B of Sheet1 per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B9").Value  = 6.380899999999993
$ws.Range("B18").Value = 6.479199999999999
$ws.Range("B20").Value = 9.309699999999996
$ws.Range("B27").Value = 6.085900000000004
$ws.Range("B69").Value = 5.422599999999993
$ws.Range("B76").Value = 5.256
$ws.Range("B82").Value = 5.589400000000003
